$d = $word.ActiveDocument

# --- Clear the entire existing body content down to a single empty paragraph ---
$d.Range(0, $d.Content.End).Delete()

# --- Build all target paragraphs in one shot, separated by paragraph marks ---
$lines = @(
    "CERTIFICADO",
    "Certificamos para os devidos fins que o colaborador:",
    "{{NOME}}",
    "Portador do CPF nº {{CPF}}",
    "Concluiu com êxito o treinamento de Norma Regulamentadora {{CURSO}}.",
    "Rio de Janeiro, {{DATA}}.",
    ""
)
$fullText = [string]::Join("`r", $lines)

$insertRange = $d.Range(0, 0)
$insertRange.Text = $fullText

# --- Apply paragraph-level + run-level formatting ---
# Para 1: "CERTIFICADO" - centered, bold, size 36 (72 half-points)
$p1 = $d.Paragraphs.Item(1)
$p1.Alignment = 1
$p1.Range.Font.Bold = 1
$p1.Range.Font.BoldBi = 1
$p1.Range.Font.Size = 36
$p1.Range.Font.SizeBi = 36

# Paras 2-6: centered, size 16 (32 half-points), not bold
for ($i = 2; $i -le 6; $i++) {
    $p = $d.Paragraphs.Item($i)
    $p.Alignment = 1
    $p.Range.Font.Size = 16
    $p.Range.Font.SizeBi = 16
}

# Para 7: trailing empty paragraph - centered, bold, size 16 (32 half-points)
$p7 = $d.Paragraphs.Item(7)
$p7.Alignment = 1
$p7.Range.Font.Bold = 1
$p7.Range.Font.BoldBi = 1
$p7.Range.Font.Size = 16
$p7.Range.Font.SizeBi = 16

# --- Page borders on the (single) section ---
$sec = $d.Sections.Item(1)
$topBorder = $sec.Borders.Item(-1)
$topBorder.LineStyle = 1
$topBorder.LineWidth = 2

Write-Host "Edit complete"
